$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "owning user" column (Q) - lets data imports carry the system/owner field.
$ws.Range("Q1").Value = "所属用户"
$ws.Range("Q2").Value = "rebuild"
$ws.Range("Q3").Value = "system"
$ws.Range("Q4").Value = "NOUSER"
$ws.Range("Q5").Value = "user@email.com"

# Q5 is an email address -> wire it up as a mailto hyperlink, like the other
# email/url columns on this sheet, then make sure the cell keeps reusing the
# workbook's existing "Hyperlink" cell style.
$ws.Hyperlinks.Add($ws.Range("Q5"), "mailto:user@email.com")
$ws.Range("Q5").Style = "超链接"

# Match the selection left behind in the authored workbook.
$ws.Range("P10").Select() | Out-Null
